# Finalized manyeyes' castle level 1
$wb = $excel.ActiveWorkbook

# --- Chests sheet ---------------------------------------------------------
$chests = $wb.Worksheets.Item("Chests")

# New chests for Manyeyes' castle level 1.
$chests.Range("A16").Value = 146
$chests.Range("B16").Value = "Manyeyes'c castle 1 (461)"
$chests.Range("C16").Value = "1x Dark Dagger"

$chests.Range("A17").Value = 147
$chests.Range("B17").Value = "Manyeyes'c castle 1 (461)"
$chests.Range("C17").Value = "1x Holy Horn, 1x Horned Helmet, 1x Scimitar, 800 Gold"

$chests.Range("A18").Value = 148
$chests.Range("B18").Value = "Manyeyes'c castle 1 (461)"
$chests.Range("C18").Value = "1x Silver Cutlery, 3 Healing Potion III, 2 Spell Potion IV, 1 Healing Potion IV, 250 Gold"

# Existing chest 143 (row 13) now also grants food.
$chests.Range("C13").Value = "1x Stamina Potion, 3x Bitter, 5x Food"

$chests.Range("C14").Select()

# --- Items sheet ------------------------------------------------------------
$items = $wb.Worksheets.Item("Items")

$items.Range("A14").Value = 415
$items.Range("B14").Value = "Dunkle Klinge / Dark Blade"
$items.Range("C14").Value = "Weapon"
$items.Range("D14").Value = "Cursed weapon from manyeyes' castle"

# --- GlobalVars sheet: move selection -----------------------------------
$globalVars = $wb.Worksheets.Item("GlobalVars")
$globalVars.Range("A29").Select()

# --- Make Items the active sheet/tab (was Monsters) ------------------------
$items.Activate()
$items.Range("E14").Select()
